# Update NATMI Wnt2-Fzd7 output with new TPM-derived values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 1.037532
$ws.Range("N2").Value = 3.112596
$ws.Range("O2").Value = 0.04166450179684251
$ws.Range("P2").Value = 0.0439159257402554
$ws.Range("Q2").Value = 0.402978812176
$ws.Range("R2").Value = 3.626809309584
$ws.Range("S2").Value = 0.04166450179684251
$ws.Range("T2").Value = 0.0439159257402554

# Row 3 (only derived-specificity columns change; M3/N3/Q3/R3 stay the same)
$ws.Range("O3").Value = 0.2561129158441639
$ws.Range("P3").Value = 0.2699524849277078
$ws.Range("S3").Value = 0.2561129158441639
$ws.Range("T3").Value = 0.2699524849277078

# Row 4
$ws.Range("M4").Value = 6.239319333333333
$ws.Range("N4").Value = 18.717958
$ws.Range("O4").Value = 0.2505543265891952
$ws.Range("P4").Value = 0.2640935262839185
$ws.Range("Q4").Value = 2.423359948159111
$ws.Range("R4").Value = 21.810239533432
$ws.Range("S4").Value = 0.2505543265891952
$ws.Range("T4").Value = 0.2640935262839185

# Row 5
$ws.Range("M5").Value = 3.8299385
$ws.Range("N5").Value = 7.659877
$ws.Range("O5").Value = 0.1538000558200097
$ws.Range("P5").Value = 0.1080739644693659
$ws.Range("Q5").Value = 1.487553219984667
$ws.Range("R5").Value = 8.925319319908001
$ws.Range("S5").Value = 0.1538000558200097
$ws.Range("T5").Value = 0.1080739644693659

# Row 6
$ws.Range("M6").Value = 7.417532333333334
$ws.Range("N6").Value = 22.252597
$ws.Range("O6").Value = 0.2978681999497886
$ws.Range("P6").Value = 0.3139640985787523
$ws.Range("Q6").Value = 2.880979448309778
$ws.Range("R6").Value = 25.92881503478801
$ws.Range("S6").Value = 0.2978681999497886
$ws.Range("T6").Value = 0.3139640985787523
